$d = $word.ActiveDocument

# 1. Bump the years-of-experience figure on the summary line.
$d.Content.Find.Execute("8.5 years", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "8.8 years", 2)

# 2. Drop the two stray blank paragraphs (an empty section-break paragraph and
#    an empty Normal-style paragraph) that trail the "12th Standard" line -
#    a leftover empty "North Eastern Trip" continuation section.
$marker = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*12th Standard*") {
        $marker = $i
        break
    }
}

if ($marker -ne $null) {
    $firstBlank = $d.Paragraphs.Item($marker + 1)
    $secondBlank = $d.Paragraphs.Item($marker + 2)
    $killRange = $d.Range($firstBlank.Range.Start, $secondBlank.Range.End)
    $killRange.Delete()
}
